$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts rows 6-10 down to 7-11)
$ws.Rows.Item(6).Insert()

# Copy formatting from row 7 (which holds the formatting that used to belong to row 6)
# down into the newly inserted blank row 6, so it matches the other interior rows.
$ws.Range("A7:C7").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's data
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Lasso Regression+normalization+ lag1"
$ws.Range("C6").Value = 97.3029684168079

# Renumber the Id column (column A) for the rows following the inserted row
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Update the active cell selection to reflect where the edit was made
$ws.Range("B6").Select() | Out-Null

$excel.CutCopyMode = $false | Out-Null
